# Add a "Source" worksheet right after Sheet1, containing a copy of
# Sheet1's data (A1:B11) — mirrors the commit "Add Source data tab to file".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet positioned immediately after Sheet1.
$source = $wb.Worksheets.Add($null, $ws1)
$source.Name = "Source"

# Copy Sheet1's used-range values into the new Source sheet.
$lastRow = $ws1.UsedRange.Rows.Count
$lastCol = $ws1.UsedRange.Columns.Count

for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $source.Cells.Item($r, $c).Value = $ws1.Cells.Item($r, $c).Value2
    }
}

# Leave Sheet1 as the active/selected tab, as in the original workbook.
$ws1.Activate()
